$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.803.84'
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('D3').Value = '2.304.46'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '304.93'
$ws.Range('E5').Value = '  +1.97%  '
$ws.Range('D6').Value = '96.69'
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').Value = '0.506'
$ws.Range('E7').Value = '  -1.11%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('D9').Value = '0.497'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').Value = '35.17'
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').Value = '18.65'
$ws.Range('E12').Value = '  +5.00%  '
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('D14').Value = '6.84'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').Value = '2.661.99'
$ws.Range('E15').Value = '  +0.52%  '
$ws.Range('D16').Value = '2.295.14'
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').Value = '0.778'
$ws.Range('E17').Value = '  +0.46%  '
$ws.Range('D18').Value = '42.736.90'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('E19').Value = '  +1.26%  '
$ws.Range('D20').Value = '0.0₃0893'
$ws.Range('E20').Value = '  -1.16%  '
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('D22').Value = '67.20'
$ws.Range('E22').Value = '  -0.85%  '
$ws.Range('D23').Value = '235.63'
$ws.Range('E23').Value = '  -2.17%  '
$ws.Range('D24').Value = '2.15'
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('E25').Value = '  +0.19%  '
$ws.Range('D26').Value = '2.41'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').Value = '24.69'
$ws.Range('E27').Value = '  -1.80%  '
$ws.Range('D28').Value = '166.02'
$ws.Range('E28').Value = '  -0.10%  '
$ws.Range('E29').Value = '  +0.84%  '
$ws.Range('D30').Value = '9.03'
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').Value = '33.07'
$ws.Range('E31').Value = '  +0.64%  '
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').Value = '17.98'
$ws.Range('E33').Value = '  +4.90%  '
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('E35').Value = '  -6.68%  '
$ws.Range('E36').Value = '  -0.76%  '
$ws.Range('D37').Value = '0.0685'
$ws.Range('E37').Value = '  +0.07%  '
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('E39').Value = '  -0.52%  '
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('E41').Value = '  -1.26%  '
$ws.Range('D42').Value = '1.998.26'
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E43').Value = '  -0.89%  '
$ws.Range('D44').Value = '10.20'
$ws.Range('E44').Value = '  +1.10%  '
$ws.Range('D45').Value = '18.16'
$ws.Range('E45').Value = '  +6.22%  '
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '2.530.00'
$ws.Range('E48').Value = '  +0.56%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').Value = '53.50'
$ws.Range('E49').Value = '  +0.81%  '
$ws.Range('E50').Value = '  +2.04%  '
$ws.Range('E51').Value = '  -1.05%  '
